$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "alpha4F"

# Tiny floating point corrections on row 13 (HKL = 11)
$ws.Range("C13").Value = 0.9896241418369076
$ws.Range("F13").Value = 0.9896241418369076
$ws.Range("H13").Value = 0.9989723178993114
$ws.Range("L13").Value = 0.9930256704571852
$ws.Range("M13").Value = 0.9947124545753637

# Tiny floating point corrections on row 15 (HKL = 13)
$ws.Range("C15").Value = 0.9749763155331567
$ws.Range("F15").Value = 0.9749763155331567

# New row 16 (HKL = 14), reuses shared string "HexGrid-60degTilt5degRes" (same as B15)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.135415534191667
$ws.Range("D16").Value = 0.6921123391403485
$ws.Range("E16").Value = 1.031329407819405
$ws.Range("F16").Value = 1.135415534191667
$ws.Range("G16").Value = 0.8370701967696368
$ws.Range("H16").Value = 1.11817816371222
$ws.Range("I16").Value = 1.066769766108889
$ws.Range("J16").Value = 0.6921123391403485
$ws.Range("K16").Value = 0.8617208734798767
$ws.Range("L16").Value = 0.9985682038357718
$ws.Range("M16").Value = 0.9801459012903612
